$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 153, pushing existing rows 153-224 down to 155-226
$ws.Rows.Item(153).Resize(2).Insert()

# New row 153 data
$ws.Range("A153").Value = 10
$ws.Range("B153").Value = "Vega Modelo de Temuco"
$ws.Range("C153").Value = "La Araucanía"
$ws.Range("D153").Value = 44466
$ws.Range("E153").Value = 9
$ws.Range("F153").Value = 100112040
$ws.Range("G153").Value = "Cilantro"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 70
$ws.Range("K153").Value = 4000
$ws.Range("L153").Value = 4000
$ws.Range("M153").Value = 4000
$ws.Range("N153").Value = "$/docena de atados (2 kilos)"
$ws.Range("O153").Value = "Provincia de Cautín"
$ws.Range("P153").Value = 2000
$ws.Range("Q153").Value = 2
$ws.Range("R153").Value = "Hortaliza"

# New row 154 data
$ws.Range("A154").Value = 10
$ws.Range("B154").Value = "Vega Modelo de Temuco"
$ws.Range("C154").Value = "La Araucanía"
$ws.Range("D154").Value = 44466
$ws.Range("E154").Value = 9
$ws.Range("F154").Value = 100112040
$ws.Range("G154").Value = "Cilantro"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 90
$ws.Range("K154").Value = 2600
$ws.Range("L154").Value = 4000
$ws.Range("M154").Value = 3378
$ws.Range("N154").Value = "$/docena de atados (2 kilos)"
$ws.Range("O154").Value = "Región Metropolitana"
$ws.Range("P154").Value = 1689
$ws.Range("Q154").Value = 2
$ws.Range("R154").Value = "Hortaliza"
